$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A held text quarter labels ("2004Q4" .. "2024Q4") as shared
# strings. Replace them with real Excel dates (Dec-31 of each year) and
# give them a custom date/time number format, matching rows 2-22.
$dates = @{
    2  = 38352   # 2004-12-31
    3  = 38717   # 2005-12-31
    4  = 39082   # 2006-12-31
    5  = 39447   # 2007-12-31
    6  = 39813   # 2008-12-31
    7  = 40178   # 2009-12-31
    8  = 40543   # 2010-12-31
    9  = 40908   # 2011-12-31
    10 = 41274   # 2012-12-31
    11 = 41639   # 2013-12-31
    12 = 42004   # 2014-12-31
    13 = 42369   # 2015-12-31
    14 = 42735   # 2016-12-31
    15 = 43100   # 2017-12-31
    16 = 43465   # 2018-12-31
    17 = 43830   # 2019-12-31
    18 = 44196   # 2020-12-31
    19 = 44561   # 2021-12-31
    20 = 44926   # 2022-12-31
    21 = 45291   # 2023-12-31
    22 = 45657   # 2024-12-31
}

foreach ($row in 2..22) {
    $cell = $ws.Range("A$row")
    $cell.Value = $dates[$row]
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
